# This workbook tracks weekly price-reporting rows for "Melón" at
# Feria Lagunitas de Puerto Montt. A new week of data (2 rows: Calameño /
# Extra and Tuna / Extra, both dated 2021-12-17) needs to be inserted right
# before the existing row 42, pushing the rest of the table down by two
# rows (the previously-last rows 143/144 end up duplicated onto the newly
# created rows 145/146, which Excel's row-insert handles automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 42 (existing rows 42..144 shift down to 44..146).
$ws.Rows.Item(42).Insert()
$ws.Rows.Item(42).Insert()

# --- New row 42: Melón / Calameño / Extra, Feria Lagunitas de Puerto Montt ---
$ws.Range("A42").Value = 4
$ws.Range("B42").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C42").Value = "Los Lagos"
$ws.Range("D42").Value = "2021-12-17"
$ws.Range("E42").Value = 10
$ws.Range("F42").Value = 100112027
$ws.Range("G42").Value = "Melón"
$ws.Range("H42").Value = "Calameño"
$ws.Range("I42").Value = "Extra"
$ws.Range("J42").Value = 200
$ws.Range("K42").Value = 15000
$ws.Range("L42").Value = 15000
$ws.Range("M42").Value = 15000
$ws.Range("N42").Value = "`$/caja 12 unidades"
$ws.Range("O42").Value = "Región de O'Higgins"
$ws.Range("P42").Value = 1250
$ws.Range("Q42").Value = 12
$ws.Range("R42").Value = "Hortaliza"

# --- New row 43: Melón / Tuna / Extra, Feria Lagunitas de Puerto Montt ---
$ws.Range("A43").Value = 4
$ws.Range("B43").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C43").Value = "Los Lagos"
$ws.Range("D43").Value = "2021-12-17"
$ws.Range("E43").Value = 10
$ws.Range("F43").Value = 100112027
$ws.Range("G43").Value = "Melón"
$ws.Range("H43").Value = "Tuna"
$ws.Range("I43").Value = "Extra"
$ws.Range("J43").Value = 250
$ws.Range("K43").Value = 15000
$ws.Range("L43").Value = 15000
$ws.Range("M43").Value = 15000
$ws.Range("N43").Value = "`$/caja 12 unidades"
$ws.Range("O43").Value = "Región de O'Higgins"
$ws.Range("P43").Value = 1250
$ws.Range("Q43").Value = 12
$ws.Range("R43").Value = "Hortaliza"
